$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are stored as text (matches source data which
# uses formatted numeric strings such as "67.560.15" / "1.00" rather than numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.560.15"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "3.265.13"
$ws.Range("E3").Value = "  -5.83%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "592.60"
$ws.Range("E5").Value = "  -3.33%  "
$ws.Range("D6").Value = "150.37"
$ws.Range("E6").Value = "  -10.61%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "3.258.78"
$ws.Range("E8").Value = "  -5.89%  "
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  -8.62%  "
$ws.Range("E10").Value = "  -10.92%  "
$ws.Range("D11").Value = "6.68"
$ws.Range("E11").Value = "  -5.97%  "
$ws.Range("D12").Value = "0.505"
$ws.Range("E12").Value = "  -10.88%  "
$ws.Range("E13").Value = "  -8.82%  "
$ws.Range("D14").Value = "38.39"
$ws.Range("E14").Value = "  -13.80%  "
$ws.Range("D15").Value = "3.792.74"
$ws.Range("D16").Value = "67.582.30"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").Value = "3.271.03"
$ws.Range("E17").Value = "  -5.66%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.114"
$ws.Range("E18").Value = "  -5.32%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "532.03"
$ws.Range("E19").Value = "  -9.30%  "
$ws.Range("D20").Value = "7.12"
$ws.Range("E20").Value = "  -13.27%  "
$ws.Range("D21").Value = "14.96"
$ws.Range("E21").Value = "  -13.23%  "
$ws.Range("D22").Value = "0.757"
$ws.Range("E22").Value = "  -11.30%  "
$ws.Range("D23").Value = "7.89"
$ws.Range("E23").Value = "  -11.80%  "
$ws.Range("D24").Value = "85.62"
$ws.Range("E24").Value = "  -10.95%  "
$ws.Range("D25").Value = "13.53"
$ws.Range("E25").Value = "  -11.32%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -10.82%  "
$ws.Range("D28").Value = "8.06"
$ws.Range("E28").Value = "  -6.97%  "
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  -11.92%  "
$ws.Range("D30").Value = "29.17"
$ws.Range("E30").Value = "  -11.66%  "
$ws.Range("D31").Value = "2.68"
$ws.Range("E31").Value = "  -5.22%  "
$ws.Range("D32").Value = "1.18"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("D33").Value = "6.62"
$ws.Range("E33").Value = "  -16.23%  "
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  -13.30%  "
$ws.Range("D35").Value = "516.25"
$ws.Range("E35").Value = "  -12.11%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "0.0442"
$ws.Range("E37").Value = "  -7.74%  "
$ws.Range("D38").Value = "53.26"
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("D39").Value = "0.0856"
$ws.Range("E39").Value = "  -11.04%  "
$ws.Range("D40").Value = "8.96"
$ws.Range("E40").Value = "  -15.54%  "
$ws.Range("E41").Value = "  -10.96%  "
$ws.Range("D42").Value = "2.77"
$ws.Range("E42").Value = "  -12.68%  "
$ws.Range("D43").Value = "2.936.69"
$ws.Range("E43").Value = "  -9.86%  "
$ws.Range("D44").Value = "0.266"
$ws.Range("E44").Value = "  -10.38%  "
$ws.Range("D45").Value = "0.0₃0589"
$ws.Range("E45").Value = "  -16.07%  "
$ws.Range("D46").Value = "2.19"
$ws.Range("E46").Value = "  -9.45%  "
$ws.Range("D47").Value = "26.70"
$ws.Range("E47").Value = "  -13.90%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  -16.49%  "
$ws.Range("D50").Value = "0.113"
$ws.Range("E50").Value = "  -10.25%  "
$ws.Range("D51").Value = "124.15"
$ws.Range("E51").Value = "  -7.17%  "
